$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "51.519.21"
Set-TextCell "E2" "  +0.83%  "
Set-TextCell "D3" "2.988.25"
Set-TextCell "E3" "  +1.53%  "
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.05%  "
Set-TextCell "D5" "382.21"
Set-TextCell "E5" "  +1.95%  "
Set-TextCell "D6" "103.66"
Set-TextCell "E6" "  +2.43%  "
Set-TextCell "D7" "0.546"
Set-TextCell "E7" "  +1.85%  "
Set-TextCell "E8" "  -0.01%  "
Set-TextCell "E9" "  +1.45%  "
Set-TextCell "D10" "36.85"
Set-TextCell "E10" "  +1.47%  "
Set-TextCell "D11" "0.137"
Set-TextCell "E11" "  -0.90%  "
Set-TextCell "E12" "  +1.27%  "
Set-TextCell "D13" "3.467.34"
Set-TextCell "E13" "  +1.84%  "
Set-TextCell "D14" "18.44"
Set-TextCell "E14" "  +2.01%  "
Set-TextCell "D15" "7.80"
Set-TextCell "E15" "  +3.14%  "
Set-TextCell "D16" "3.006.47"
Set-TextCell "E16" "  +1.74%  "
Set-TextCell "D17" "11.12"
Set-TextCell "E17" "  +4.40%  "
Set-TextCell "D18" "0.998"
Set-TextCell "E18" "  +0.67%  "
Set-TextCell "D19" "51.517.55"
Set-TextCell "E19" "  +1.06%  "
Set-TextCell "D20" "3.08"
Set-TextCell "E20" "  -0.80%  "
Set-TextCell "D21" "12.62"
Set-TextCell "E21" "  +1.47%  "
Set-TextCell "D22" "0.0₃0964"
Set-TextCell "E22" "  +0.76%  "
Set-TextCell "D23" "70.57"
Set-TextCell "E23" "  +2.75%  "
Set-TextCell "D24" "267.81"
Set-TextCell "E24" "  +0.71%  "
Set-TextCell "E25" "  +2.75%  "
Set-TextCell "D26" "7.88"
Set-TextCell "E26" "  -2.98%  "
Set-TextCell "D27" "7.42"
Set-TextCell "E27" "  -2.74%  "
Set-TextCell "E28" "  -0.11%  "
Set-TextCell "B29" "EthereumClassic"
Set-TextCell "C29" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D29" "26.05"
Set-TextCell "E29" "  +1.69%  "
Set-TextCell "B30" "Kaspa"
Set-TextCell "C30" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D30" "0.166"
Set-TextCell "E30" "  +1.38%  "
Set-TextCell "E31" "  -0.96%  "
Set-TextCell "D32" "10.35"
Set-TextCell "E32" "  +3.18%  "
Set-TextCell "D33" "34.66"
Set-TextCell "E33" "  +4.03%  "
Set-TextCell "E35" "  +1.14%  "
Set-TextCell "D36" "0.0441"
Set-TextCell "E36" "  -0.26%  "
Set-TextCell "E37" "  +0.07%  "
Set-TextCell "D38" "3.25"
Set-TextCell "E38" "  +2.76%  "
Set-TextCell "D39" "16.84"
Set-TextCell "E39" "  +3.58%  "
Set-TextCell "B40" "Stellar"
Set-TextCell "C40" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D40" "0.116"
Set-TextCell "E40" "  +1.38%  "
Set-TextCell "B41" "Stacks"
Set-TextCell "C41" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D41" "2.56"
Set-TextCell "E41" "  +3.11%  "
Set-TextCell "D42" "1.84"
Set-TextCell "E42" "  +2.88%  "
Set-TextCell "D43" "124.65"
Set-TextCell "E43" "  +3.54%  "
Set-TextCell "E44" "  +9.56%  "
Set-TextCell "D45" "21.52"
Set-TextCell "E45" "  +0.39%  "
Set-TextCell "D46" "2.03"
Set-TextCell "E46" "  +0.09%  "
Set-TextCell "D47" "2.39"
Set-TextCell "E47" "  +3.30%  "
Set-TextCell "D48" "0.270"
Set-TextCell "E48" "  -0.87%  "
Set-TextCell "D49" "2.043.47"
Set-TextCell "E49" "  +2.33%  "
Set-TextCell "D50" "0.0332"
Set-TextCell "E50" "  +2.46%  "
Set-TextCell "D51" "0.539"
Set-TextCell "E51" "  +16.51%  "
